$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The new TPM run drops the "Neutrophils" target-cluster rows entirely (rows 6 and 12 of the
# original 13-row table). Delete the lower row first so the upper row index is unaffected.
$ws.Rows.Item(12).Delete()
$ws.Rows.Item(6).Delete()

# The surviving 10 rows (now rows 2-11) all get refreshed numeric columns (E:T) from the new TPM data.
# Row 2: FAPs -> ECs
$ws.Range("E2").Value = 3.0
$ws.Range("F2").Value = 1.0
$ws.Range("G2").Value = 7.514794999999999
$ws.Range("H2").Value = 22.544385
$ws.Range("I2").Value = 0.977669497583861
$ws.Range("J2").Value = 0.977669497583861
$ws.Range("K2").Value = 3.0
$ws.Range("L2").Value = 1.0
$ws.Range("M2").Value = 12.13731566666667
$ws.Range("N2").Value = 36.411947
$ws.Range("O2").Value = 0.2052045499593414
$ws.Range("P2").Value = 0.2052045499593414
$ws.Range("Q2").Value = 91.20943908528832
$ws.Range("R2").Value = 820.8849517675949
$ws.Range("S2").Value = 0.2006222292606716
$ws.Range("T2").Value = 0.2006222292606716

# Row 3: FAPs -> FAPs
$ws.Range("E3").Value = 3.0
$ws.Range("F3").Value = 1.0
$ws.Range("G3").Value = 7.514794999999999
$ws.Range("H3").Value = 22.544385
$ws.Range("I3").Value = 0.977669497583861
$ws.Range("J3").Value = 0.977669497583861
$ws.Range("K3").Value = 3.0
$ws.Range("L3").Value = 1.0
$ws.Range("M3").Value = 35.95218066666666
$ws.Range("N3").Value = 107.856542
$ws.Range("O3").Value = 0.6078404201038963
$ws.Range("P3").Value = 0.6078404201038962
$ws.Range("Q3").Value = 270.1732675129633
$ws.Range("R3").Value = 2431.559407616669
$ws.Range("S3").Value = 0.5942670381341393
$ws.Range("T3").Value = 0.5942670381341392

# Row 4: FAPs -> Inflammatory-Mac
$ws.Range("E4").Value = 3.0
$ws.Range("F4").Value = 1.0
$ws.Range("G4").Value = 7.514794999999999
$ws.Range("H4").Value = 22.544385
$ws.Range("I4").Value = 0.977669497583861
$ws.Range("J4").Value = 0.977669497583861
$ws.Range("K4").Value = 2.0
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 0.101426
$ws.Range("N4").Value = 0.304278
$ws.Range("O4").Value = 0.001714800640914052
$ws.Range("P4").Value = 0.001714800640914052
$ws.Range("Q4").Value = 0.76219559767
$ws.Range("R4").Value = 6.859760379029999
$ws.Range("S4").Value = 0.001676508281058924
$ws.Range("T4").Value = 0.001676508281058924

# Row 5: FAPs -> MuSCs
$ws.Range("E5").Value = 3.0
$ws.Range("F5").Value = 1.0
$ws.Range("G5").Value = 7.514794999999999
$ws.Range("H5").Value = 22.544385
$ws.Range("I5").Value = 0.977669497583861
$ws.Range("J5").Value = 0.977669497583861
$ws.Range("K5").Value = 3.0
$ws.Range("L5").Value = 1.0
$ws.Range("M5").Value = 10.82580366666667
$ws.Range("N5").Value = 32.477411
$ws.Range("O5").Value = 0.1830309296039447
$ws.Range("P5").Value = 0.1830309296039447
$ws.Range("Q5").Value = 81.35369526524833
$ws.Range("R5").Value = 732.183257387235
$ws.Range("S5").Value = 0.1789437569881956
$ws.Range("T5").Value = 0.1789437569881956

# Row 6: FAPs -> Resolving-Mac
$ws.Range("E6").Value = 3.0
$ws.Range("F6").Value = 1.0
$ws.Range("G6").Value = 7.514794999999999
$ws.Range("H6").Value = 22.544385
$ws.Range("I6").Value = 0.977669497583861
$ws.Range("J6").Value = 0.977669497583861
$ws.Range("K6").Value = 3.0
$ws.Range("L6").Value = 1.0
$ws.Range("M6").Value = 0.1306743333333333
$ws.Range("N6").Value = 0.392023
$ws.Range("O6").Value = 0.002209299691903619
$ws.Range("P6").Value = 0.002209299691903619
$ws.Range("Q6").Value = 0.9819908267616666
$ws.Range("R6").Value = 8.837917440855
$ws.Range("S6").Value = 0.00215996491979559
$ws.Range("T6").Value = 0.00215996491979559

# Row 7: MuSCs -> ECs
$ws.Range("E7").Value = 2.0
$ws.Range("F7").Value = 0.6666666666666666
$ws.Range("G7").Value = 0.171642
$ws.Range("H7").Value = 0.514926
$ws.Range("I7").Value = 0.02233050241613897
$ws.Range("J7").Value = 0.02233050241613898
$ws.Range("K7").Value = 3.0
$ws.Range("L7").Value = 1.0
$ws.Range("M7").Value = 12.13731566666667
$ws.Range("N7").Value = 36.411947
$ws.Range("O7").Value = 0.2052045499593414
$ws.Range("P7").Value = 0.2052045499593414
$ws.Range("Q7").Value = 2.083273135658
$ws.Range("R7").Value = 18.749458220922
$ws.Range("S7").Value = 0.004582320698669784
$ws.Range("T7").Value = 0.004582320698669783

# Row 8: MuSCs -> FAPs
$ws.Range("E8").Value = 2.0
$ws.Range("F8").Value = 0.6666666666666666
$ws.Range("G8").Value = 0.171642
$ws.Range("H8").Value = 0.514926
$ws.Range("I8").Value = 0.02233050241613897
$ws.Range("J8").Value = 0.02233050241613898
$ws.Range("K8").Value = 3.0
$ws.Range("L8").Value = 1.0
$ws.Range("M8").Value = 35.95218066666666
$ws.Range("N8").Value = 107.856542
$ws.Range("O8").Value = 0.6078404201038963
$ws.Range("P8").Value = 0.6078404201038962
$ws.Range("Q8").Value = 6.170904193987999
$ws.Range("R8").Value = 55.53813774589199
$ws.Range("S8").Value = 0.01357338196975698
$ws.Range("T8").Value = 0.01357338196975698

# Row 9: MuSCs -> Inflammatory-Mac
$ws.Range("E9").Value = 2.0
$ws.Range("F9").Value = 0.6666666666666666
$ws.Range("G9").Value = 0.171642
$ws.Range("H9").Value = 0.514926
$ws.Range("I9").Value = 0.02233050241613897
$ws.Range("J9").Value = 0.02233050241613898
$ws.Range("K9").Value = 2.0
$ws.Range("L9").Value = 0.6666666666666666
$ws.Range("M9").Value = 0.101426
$ws.Range("N9").Value = 0.304278
$ws.Range("O9").Value = 0.001714800640914052
$ws.Range("P9").Value = 0.001714800640914052
$ws.Range("Q9").Value = 0.017408961492
$ws.Range("R9").Value = 0.156680653428
$ws.Range("S9").Value = 0.0000382923598551279
$ws.Range("T9").Value = 0.00003829235985512789

# Row 10: MuSCs -> MuSCs
$ws.Range("E10").Value = 2.0
$ws.Range("F10").Value = 0.6666666666666666
$ws.Range("G10").Value = 0.171642
$ws.Range("H10").Value = 0.514926
$ws.Range("I10").Value = 0.02233050241613897
$ws.Range("J10").Value = 0.02233050241613898
$ws.Range("K10").Value = 3.0
$ws.Range("L10").Value = 1.0
$ws.Range("M10").Value = 10.82580366666667
$ws.Range("N10").Value = 32.477411
$ws.Range("O10").Value = 0.1830309296039447
$ws.Range("P10").Value = 0.1830309296039447
$ws.Range("Q10").Value = 1.858162592954
$ws.Range("R10").Value = 16.723463336586
$ws.Range("S10").Value = 0.00408717261574905
$ws.Range("T10").Value = 0.00408717261574905

# Row 11: MuSCs -> Resolving-Mac
$ws.Range("E11").Value = 2.0
$ws.Range("F11").Value = 0.6666666666666666
$ws.Range("G11").Value = 0.171642
$ws.Range("H11").Value = 0.514926
$ws.Range("I11").Value = 0.02233050241613897
$ws.Range("J11").Value = 0.02233050241613898
$ws.Range("K11").Value = 3.0
$ws.Range("L11").Value = 1.0
$ws.Range("M11").Value = 0.1306743333333333
$ws.Range("N11").Value = 0.392023
$ws.Range("O11").Value = 0.002209299691903619
$ws.Range("P11").Value = 0.002209299691903619
$ws.Range("Q11").Value = 0.022429203922
$ws.Range("R11").Value = 0.201862835298
$ws.Range("S11").Value = 0.00004933477210802885
$ws.Range("T11").Value = 0.00004933477210802885

